# Training Results.xlsx - incorporate LabelEncoding results
# - Rename existing "No Dum" sheet to "No PCA or Sampling - No Dum" (matches the
#   naming convention already used by the "- Dum" sheets) and refresh its values
#   (it now holds the "No PCA / No Sampling" results instead of the old placeholder).
# - Add three sibling "- No Dum" sheets (No PCA & Sampling, PCA & No Sampling,
#   PCA & Sampling), mirroring sheets 1-4, each populated with new results.
# - Refresh the metric values on the four existing "- Dum" sheets.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1. Update metric values on the four existing "- Dum" sheets (sheet1..sheet4)
# ---------------------------------------------------------------------------

$ws1 = $wb.Worksheets.Item("No PCA or Sampling - Dum")
$ws1.Range("B2").Value = 0.7328008386888083
$ws1.Range("C2").Value = 0.7182267694502146
$ws1.Range("D2").Value = 0.6666712784158447
$ws1.Range("B3").Value = 0.754846626109054
$ws1.Range("C3").Value = 0.7296295728980202
$ws1.Range("D3").Value = 0.6998728061898157
$ws1.Range("B4").Value = 0.7492060273464546
$ws1.Range("C4").Value = 0.7130707888411264
$ws1.Range("D4").Value = 0.6789985597219869

$ws2 = $wb.Worksheets.Item("No PCA & Sampling - Dum")
$ws2.Range("B2").Value = 0.8991316141252031
$ws2.Range("C2").Value = 0.6062146957609581
$ws2.Range("D2").Value = 0.6681878164082138
$ws2.Range("B3").Value = 0.8977984332326566
$ws2.Range("C3").Value = 0.7087584538846892
$ws2.Range("D3").Value = 0.752196484509269
$ws2.Range("B4").Value = 0.9081247104456669
$ws2.Range("C4").Value = 0.710620984889252
$ws2.Range("D4").Value = 0.7410838929053289

$ws3 = $wb.Worksheets.Item("PCA & No Sampling - Dum")
$ws3.Range("C2").Value = 0.5241990365991013
$ws3.Range("D2").Value = 0.5057262920367432
$ws3.Range("C3").Value = 0.5439637384008177
$ws3.Range("D3").Value = 0.5041115301326291
$ws3.Range("C4").Value = 0.5295187669020757
$ws3.Range("D4").Value = 0.5061023341830998

$ws4 = $wb.Worksheets.Item("PCA & Sampling - Dum")
$ws4.Range("B2").Value = 0.5703971105767545
$ws4.Range("C2").Value = 0.5171297512383519
$ws4.Range("D2").Value = 0.5082944452455465
$ws4.Range("B3").Value = 0.5761295204104899
$ws4.Range("C3").Value = 0.5305166704308398
$ws4.Range("D3").Value = 0.5101064094142234
$ws4.Range("B4").Value = 0.5734577039354652
$ws4.Range("C4").Value = 0.5172688215141598
$ws4.Range("D4").Value = 0.5051671342469406

# ---------------------------------------------------------------------------
# 2. Rename the old "No Dum" sheet and refresh/extend its values to 3 rows
# ---------------------------------------------------------------------------

$ws5 = $wb.Worksheets.Item("No PCA & Sampling - No Dum")
$ws5.Name = "No PCA or Sampling - No Dum"

# Extend with rows 3 and 4 (copy index-column + header styling down/along)
$ws5.Range("A2").Copy($ws5.Range("A3"))
$ws5.Range("A2").Copy($ws5.Range("A4"))

$ws5.Range("A2").Value = 0
$ws5.Range("B2").Value = 0.7360908410661325
$ws5.Range("C2").Value = 0.7358185594974103
$ws5.Range("D2").Value = 0.7801402036704326

$ws5.Range("A3").Value = 1
$ws5.Range("B3").Value = 0.7623299701632741
$ws5.Range("C3").Value = 0.8183232386063142
$ws5.Range("D3").Value = 0.7953141788673466

$ws5.Range("A4").Value = 2
$ws5.Range("B4").Value = 0.7421561252543215
$ws5.Range("C4").Value = 0.8075634211136422
$ws5.Range("D4").Value = 0.7883506417993209

# ---------------------------------------------------------------------------
# 3. Add the three remaining "- No Dum" sheets, mirroring the "- Dum" layout
# ---------------------------------------------------------------------------

function Add-ResultsSheet {
    param(
        [string]$Name,
        [double[][]]$Rows
    )

    $lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
    $newSheet = $wb.Worksheets.Add($null, $lastSheet)
    $newSheet.Name = $Name

    # Bring over header + index-column formatting from the template sheet
    $ws1.Range("B1:D1").Copy($newSheet.Range("B1"))
    $ws1.Range("A2:A4").Copy($newSheet.Range("A2"))

    for ($i = 0; $i -lt $Rows.Length; $i++) {
        $r = $i + 2
        $newSheet.Range("A$r").Value = $Rows[$i][0]
        $newSheet.Range("B$r").Value = $Rows[$i][1]
        $newSheet.Range("C$r").Value = $Rows[$i][2]
        $newSheet.Range("D$r").Value = $Rows[$i][3]
    }

    return $newSheet
}

$rows6 = @(
    @(0, 0.8766352093409355, 0.7706789599011655, 0.7904237872723517),
    @(1, 0.8800602079005679, 0.7868627291667586, 0.8126572908925519),
    @(2, 0.8961957415257527, 0.7900933246269107, 0.8209197092470765)
)
$ws6 = Add-ResultsSheet "No PCA & Sampling - No Dum" $rows6

$rows7 = @(
    @(0, 0.5, 0.5202473001781045, 0.5061405310396702),
    @(1, 0.5, 0.5258897135845804, 0.5058532961461333),
    @(2, 0.5, 0.5455169761456464, 0.5076209164321281)
)
$ws7 = Add-ResultsSheet "PCA & No Sampling - No Dum" $rows7

$rows8 = @(
    @(0, 0.5447909903309157, 0.5169213066709276, 0.5041336212104087),
    @(1, 0.5561560588716121, 0.519904074772212, 0.5082789193179668),
    @(2, 0.5541294770781681, 0.5223841402569087, 0.509912437903126)
)
$ws8 = Add-ResultsSheet "PCA & Sampling - No Dum" $rows8
